# Update the "Förändrad" (Changed) date column (C) for rows 2-6
# from 2023-09-01 (serial 45170) to 2023-09-05 (serial 45174).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDate = Get-Date -Year 2023 -Month 9 -Day 5 -Hour 0 -Minute 0 -Second 0

for ($row = 2; $row -le 6; $row++) {
    $ws.Cells.Item($row, 3).Value = $newDate
}
